# Updated via Streamlit Approval System
# Clear the approval/ledger columns (AI:AO -> APPROVAL_1, APPROVAL_2,
# COST_CENTER, LEDGER_NAME, LEDGER_UNDER, TO, BY) for all data rows
# (rows 2-25), resetting them back to blank/empty text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("AI2:AO25")

# Assign a lone apostrophe so Excel commits the edit as an (empty) TEXT
# value rather than normalizing the blank straight through to a generic
# "no value" cell - this keeps the cells' text type, matching how the
# sheet already represents its other blank string columns.
$rng.Value = "'"

# Re-apply the default "Normal" cell style so no stray quote-prefix
# formatting is left behind on the now-empty cells.
$rng.Style = "Normal"
